$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Widen column B to fit the new, longer parameter names
$ws.Columns.Item(2).ColumnWidth = 27

# New parameter rows - entered in the same order as the original authoring
# session so the shared-strings table comes out in the same sequence.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "p_houseOtherElectricityDemandPeak_kW"
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = "kW"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "p_storeOtherElectricityDemandPeak_kW"
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = "kW"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "p_officeOtherElectricityDemandPeak_kW"
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = "kW"

$ws.Range("E3:E5").Value = "Peak power demand for lighting, devices. To be scaled by profile"

$ws.Range("A6").Value = 4
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = "kW"
$ws.Range("E6").Value = "House GasBurner peak thermal delivery in kW"
$ws.Range("B6").Value = "p_houseGasBurnerThermalCapacity_kW"

# Update selection to match the saved view state
$ws.Range("B7").Select()
